$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update header cell A1 from "VISUALIZATION" to "Type"
$ws.Range("A1").Value = "Type"

# Update the selected/active cell to A12 (matches new sheetView selection)
$ws.Range("A12").Select()
